$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Cadastro"

# Insert a new column for "Pais" (Country) before the existing "Cidade" column (col H / 8)
$ws.Columns.Item(8).Insert()

# Match the width used by the neighboring "Numero de Telefone" column (col G)
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(7).ColumnWidth

# Reformat phone numbers as text strings with the local formatting
$ws.Range("G2").Value = "(11)970510060"
$ws.Range("G3").Value = "(11)970510070"

# Fill in header + data for the new column
$ws.Range("H1").Value = "Pais"
$ws.Range("H2").Value = "Brazil"
$ws.Range("H3").Value = "Brazil"

# Leave the selection on the newly added cell, like the author did
$ws.Range("H3").Select()
